$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the bullet / balancing values
$ws.Range("B59").Value = 300
$ws.Range("B63").Value = 700
$ws.Range("B119").Value = 2
$ws.Range("B138").Value = 10
$ws.Range("B139").Value = 36
$ws.Range("B151").Value = 1000
$ws.Range("B163").Value = 3

# Update the sheet view (scroll position & selection)
$excel.ActiveWindow.ScrollRow = 100
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D117").Select()
